$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - RandomForestRegressor (name unchanged), update metric values
$ws.Range("B3").Value = 38770156136267.39
$ws.Range("C3").Value = 45928845474583.34
$ws.Range("D3").Value = 15556388594436.73

# Row 4 - rename model and update metric values
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 28061926624961.07
$ws.Range("C4").Value = 16200905954044.45
$ws.Range("D4").Value = 11589858195980.47

# Row 5 - rename model and update metric values
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 125166664809328
$ws.Range("C5").Value = 107286479372162.4
$ws.Range("D5").Value = 127533652386852.7
